$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns E,F,G,H,M,N,O,P,Q,R,S,T
$data = @{
    2 = @{ E=2; F=0.6666666666666666; G=0.1618313333333333; H=0.485494; M=61.89239633333334; N=185.677189; O=0.3758067454097886; P=0.3758067454097886; Q=10.01612902181845; R=90.145161196366; S=0.3758067454097886; T=0.3758067454097886 }
    3 = @{ E=2; F=0.6666666666666666; G=0.1618313333333333; H=0.485494; O=0.00645640217309452; P=0.006456402173094521; Q=0.1720782236411111; R=1.54870401277; S=0.00645640217309452; T=0.006456402173094521 }
    4 = @{ E=2; F=0.6666666666666666; G=0.1618313333333333; H=0.485494; M=11.19030466666667; N=33.570914; O=0.06794682749517447; P=0.06794682749517447; Q=1.810941924612889; R=16.298477321516; S=0.06794682749517447; T=0.06794682749517447 }
    5 = @{ E=2; F=0.6666666666666666; G=0.1618313333333333; H=0.485494; M=74.55619899999999; N=223.668597; O=0.4527005602661487; P=0.4527005602661487; Q=12.06552909243533; R=108.589761831918; S=0.4527005602661487; T=0.4527005602661487 }
    6 = @{ E=2; F=0.6666666666666666; G=0.1618313333333333; H=0.485494; M=3.330078; N=9.990233999999999; O=0.02022002457944478; P=0.02022002457944478; Q=0.538910962844; R=4.850198665595999; S=0.02022002457944478; T=0.02022002457944478 }
    7 = @{ E=2; F=0.6666666666666666; G=0.1618313333333333; H=0.485494; M=12.65978833333333; N=37.979365; O=0.07686944007634902; P=0.07686944007634904; Q=2.048750425701111; R=18.43875383131; S=0.07686944007634902; T=0.07686944007634904 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
